# Fill in the "My Info" / "Lineup" scorecard form on the Lineup_Template
# sheet with a submitted entry, then leave the selection on B18 (matching
# the author filling the form top to bottom and finishing on the last cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lineup_Template")

$ws.Range("B2").Value = "Albert O"
$ws.Range("B3").Value = 1
$ws.Range("B5").Value = "Alfred L"
$ws.Range("B6").Value = "Mike K"
$ws.Range("B7").Value = "Leah M"
$ws.Range("B8").Value = "Adam A"
$ws.Range("B9").Value = "John J"
$ws.Range("B10").Value = "Brooks K"
$ws.Range("B11").Value = "Brooks K"
$ws.Range("B12").Value = "Scott S"
$ws.Range("B13").Value = "Bryson D"
$ws.Range("B14").Value = "Christy J"
$ws.Range("B15").Value = "Albert O"
$ws.Range("B16").Value = "Chet H"
$ws.Range("B17").Value = "Christy J"
$ws.Range("B18").Value = "Houston"

$ws.Range("B18").Select()
